$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table cells per the latest scrape.
# Some "Price" values look numeric (e.g. "301.70"); assigning them directly
# via .Value would make Excel auto-convert them to numbers (dropping the
# trailing zero / changing formatting). Force text entry, then clear the
# temporary format so the cell style matches the rest of the sheet (General).
$ws.Range("D2").Value = "44.193.74"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.219.42"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.89"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "2.554.99"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "2.288.54"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.799"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.10"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").Value = "43.910.31"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "0.0₃0901"
$ws.Range("E19").Value = "  -6.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.39"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.38"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.04"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0752"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.92"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -4.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.71"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.29"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0286"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "1.741.70"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "77.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "66.48"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("E50").Value = "  -6.86%  "
$ws.Range("D51").Value = "2.436.28"
$ws.Range("E51").Value = "  -1.06%  "
